# Saldo.xlsx update:
#  - Row 5 (account 004381180 / HFR / 50543.01) is replaced by two new rows:
#      002973105 / DARLAN   / 54850.53
#      005009947 / VERANICE / 47959.4
#  - 58 other account rows (various small balances) are removed entirely.
#
# All row numbers below refer to 1-based Excel rows in the ORIGINAL sheet
# layout (row 1 = header "Conta/Nome/Saldo", row 2 = first data row, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Delete the removed-account rows, working from the bottom of the
#        sheet upward (in contiguous blocks) so earlier row numbers stay
#        valid while later deletions happen. ---
$ws.Rows("63:71").Delete()
$ws.Rows("19:61").Delete()
$ws.Rows("14:17").Delete()
$ws.Rows("12:12").Delete()
$ws.Rows("9:9").Delete()

# --- 2) Replace row 5 (004381180 / HFR / 50543.01) with two rows. ---
# Insert a fresh row right after row 5, pushing everything below down one.
$ws.Rows(6).Insert()

# Row 5 becomes the first new record. The "Conta" account numbers are
# all-digit strings with significant leading zeros (e.g. "002973105"), so
# force the cell to text format before assigning, then clear the format
# override back off again (matching the rest of the sheet, which has no
# explicit per-cell style on the data rows - only the inlineStr cell type
# keeps the leading zeros intact).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "002973105"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "DARLAN"
$ws.Range("C5").Value = 54850.53

# Row 6 (the newly inserted, previously empty row) becomes the second
# new record.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "005009947"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = "VERANICE"
$ws.Range("C6").Value = 47959.4

Write-Output "Edit complete."
